# Updated cryptos list values (Price / Volume(1h)), matching the target diff.
# All target cells are plain text in the workbook (inline strings), so force
# the cell number format to Text before assigning, to stop Excel's COM layer
# from auto-coercing numeric-looking strings (e.g. "581.32", "0.0524") into
# real numbers / scientific notation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.472.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.08%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.448.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.16%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.32"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.85%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.64"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.67%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.529"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.55%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.447.49"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.14%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.63%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.16%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.21"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.33%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.347"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.65%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.52"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.02%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.95%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.859.91"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.27%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.251.10"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.41%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.430.17"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.36%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.03"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.72%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.14"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.83%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "333.25"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.75%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.33%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.97"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -7.24%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.27%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.94"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.17%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "636.52"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.14%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.17"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.55%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0963"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -8.58%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.556.16"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.66%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.02%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.44"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -6.33%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.09"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.62%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.53%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.89"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.46%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.99"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.76%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.19%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.45"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.47%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.377"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.01%  "

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.44"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.34%  "

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "148.87"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.93%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.29"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.79%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.74"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.78%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.71"

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.01%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.51"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -8.93%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "144.50"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.49%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.67"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.21%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0524"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.76%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.599"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.30%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.81"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -7.85%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0234"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.77%  "
